$d = $word.ActiveDocument
$range = $d.Content
$found = $range.Find.Execute("and also", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    Write-Host "WordOpenXML:"
    Write-Host $range.WordOpenXML
}
